$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A85:D85").Copy()
$ws.Range("A86:D86").PasteSpecial(-4122)

$ws.Range("A86").Value = "Bassem Nabil"
$ws.Range("B86").Value = "01224446379"
$ws.Range("C86").Value = "#185"
$ws.Range("D86").Value = "#85"

$ws.Range("A86:D86").EntireRow.AutoFit()

$ws.Range("C86").Select()
